$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 254.19
$ws.Range("C2").Value = 216.06
$ws.Range("D2").Value = 292.32

$ws.Range("B3").Value = 395.7
$ws.Range("C3").Value = 336.34
$ws.Range("D3").Value = 455.05

$ws.Range("B4").Value = 287.32
$ws.Range("C4").Value = 244.23
$ws.Range("D4").Value = 330.42
